$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits right after "ROSA QUILINDO".
# It needs to end up on the very last paragraph once the new content below
# has been appended, so pull it out now and re-add it at the end later.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Build the replacement text: the existing "ROSA QUILINDO" line, followed by
# six blank centered paragraphs, then the two new institution lines, then a
# trailing blank paragraph that will receive the bookmark back. Each blank
# line is tagged with a unique placeholder (B1..B6 / BB) so it can be
# targeted individually afterwards; Find/Replace on the lone, still-unique
# "ROSA QUILINDO" anchor lets every new paragraph inherit its run formatting
# (centered, es-ES) in one shot.
$sena = "SERVICIO NACIONAL DE APRENDIZAJE " + [char]8220 + "SENA" + [char]8221
$replacement = "ROSA QUILINDO^pB1^pB2^pB3^pB4^pB5^pB6^p" + $sena + "^pCAUCA^pBB"

$find_r = $d.Content
$find_r.Find.Execute("ROSA QUILINDO", $false, $false, $false, $false, $false, `
    $true, 1, $false, $replacement, 2) | Out-Null

# Strip each placeholder back out, leaving genuinely empty paragraphs behind.
foreach ($mark in @("B1", "B2", "B3", "B4", "B5", "B6", "BB")) {
    $fr = $d.Content
    $fr.Find.Execute($mark, $false, $false, $false, $false, $false, `
        $true, 1, $false, "", 2) | Out-Null
}

# Re-attach the bookmark to the new final (empty) paragraph.
$lastPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
